# Auto-generated edit script: update "想去人数" (column F) counts
# per sheet, matching the target diff (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1012
$ws.Range("F5").Value = 973
$ws.Range("F6").Value = 1550
$ws.Range("F7").Value = 39972
$ws.Range("F8").Value = 5
$ws.Range("F10").Value = 8506
$ws.Range("F12").Value = 558
$ws.Range("F13").Value = 716
$ws.Range("F14").Value = 570
$ws.Range("F15").Value = 119
$ws.Range("F16").Value = 217
$ws.Range("F17").Value = 698
$ws.Range("F19").Value = 95
$ws.Range("F20").Value = 574
$ws.Range("F21").Value = 234
$ws.Range("F22").Value = 1099
$ws.Range("F24").Value = 369
$ws.Range("F25").Value = 569
$ws.Range("F26").Value = 394
$ws.Range("F27").Value = 581
$ws.Range("F28").Value = 595
$ws.Range("F30").Value = 25
$ws.Range("F31").Value = 20
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 403
$ws.Range("F34").Value = 7
$ws.Range("F35").Value = 149
$ws.Range("F36").Value = 842
$ws.Range("F37").Value = 376
$ws.Range("F38").Value = 25
$ws.Range("F42").Value = 1034
$ws.Range("F43").Value = 218
$ws.Range("F44").Value = 1064
$ws.Range("F45").Value = 342
$ws.Range("F47").Value = 8

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 198
$ws.Range("F5").Value = 331
$ws.Range("F6").Value = 4389
$ws.Range("F8").Value = 309
$ws.Range("F12").Value = 98
$ws.Range("F17").Value = 74
$ws.Range("F20").Value = 4367

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1815
$ws.Range("F4").Value = 361
$ws.Range("F5").Value = 169

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1815
$ws.Range("F4").Value = 361
$ws.Range("F7").Value = 973
$ws.Range("F8").Value = 1551
$ws.Range("F9").Value = 39972
$ws.Range("F10").Value = 5
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = 198
$ws.Range("F13").Value = 331
$ws.Range("F14").Value = 309
$ws.Range("F17").Value = 558
$ws.Range("F19").Value = 169
$ws.Range("F20").Value = 716
$ws.Range("F21").Value = 98
$ws.Range("F22").Value = 98
$ws.Range("F23").Value = 119
$ws.Range("F24").Value = 217
$ws.Range("F25").Value = 698
$ws.Range("F27").Value = 95
$ws.Range("F28").Value = 234
$ws.Range("F29").Value = 1099
$ws.Range("F30").Value = 394
$ws.Range("F31").Value = 581
$ws.Range("F32").Value = 595
$ws.Range("F34").Value = 20
$ws.Range("F35").Value = 7
$ws.Range("F37").Value = 403
$ws.Range("F38").Value = 7
$ws.Range("F39").Value = 149
$ws.Range("F40").Value = 842
$ws.Range("F41").Value = 376
$ws.Range("F42").Value = 25
$ws.Range("F46").Value = 218
$ws.Range("F47").Value = 1064
$ws.Range("F48").Value = 342

